$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header text in D1 from "Created date & time" to "Last modified date & time"
$ws.Range("D1").Value = "Last modified date & time"

# Widen column D to fit the new, longer header text
# (24.5 is the ColumnWidth value that rounds to the target stored width of
# 25.33203125 / closest reachable value through Excel's pixel-quantized
# ColumnWidth property)
$ws.Columns.Item(4).ColumnWidth = 24.5

# Update the active cell selection to H16 (matches the recorded UI state)
$ws.Range("H16").Select()
